$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-03-05 Wednesday"; new = "2025-03-06 Thursday"},
    @{old = "391×2="; new = "560×7="},
    @{old = "506×8="; new = "402×9="},
    @{old = "676×7="; new = "458×9="},
    @{old = "836×9="; new = "378×5="},
    @{old = "846×7="; new = "142×5="},
    @{old = "836×5="; new = "197×9="},
    @{old = "671×6="; new = "822×8="},
    @{old = "679×5="; new = "890×9="},
    @{old = "728×8="; new = "984×3="},
    @{old = "532×9="; new = "726×5="},
    @{old = "670×2="; new = "414×5="},
    @{old = "562×4="; new = "701×2="},
    @{old = "240×8="; new = "742×4="},
    @{old = "119×4="; new = "954×5="},
    @{old = "443×2="; new = "143×9="},
    @{old = "531×5="; new = "713×8="},
    @{old = "909×3="; new = "977×8="},
    @{old = "382×5="; new = "767×8="},
    @{old = "782×3="; new = "146×4="},
    @{old = "649×9="; new = "293×9="},
    @{old = "491×2="; new = "954×3="},
    @{old = "718×2="; new = "533×8="},
    @{old = "344×3="; new = "922×5="},
    @{old = "815×4="; new = "973×3="},
    @{old = "142×7="; new = "555×4="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
